$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.056.45'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.812.21'
$ws.Range("E3").Value = '  -0.82%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  -0.09%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '698.99'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.53%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.00'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '3.807.88'
$ws.Range("E8").Value = '  +0.05%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  -0.42%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.52'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +3.02%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("E13").Value = '  -0.82%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.00'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").Value = '3.817.80'
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = '71.089.50'
$ws.Range("E17").Value = '  +0.21%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.49'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  -0.53%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '512.51'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +4.21%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.68'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -0.45%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.10%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.99'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").Value = '3.966.89'
$ws.Range("E26").Value = '  -0.80%  '
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("E29").Value = '  +0.07%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -3.91%  '
$ws.Range("E31").Value = '  -4.40%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.42'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -0.86%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.25'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("E35").Value = '  -4.74%  '
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D37").Value = '3.776.55'
$ws.Range("E37").Value = '  -0.61%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("E40").Value = '  +0.90%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.00'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("E42").Value = '  -1.00%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.31'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("E44").Value = '  -0.01%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '173.07'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +5.65%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("E47").Value = '  +0.00%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.41'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +1.44%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '428.86'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +5.01%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.66'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.294'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -1.25%  '
